# Mississippi.xlsx update:
#  - "Data" sheet header cells get their shorter labels
#      A1: "Interest Rates (x)"  -> "Interest Rates"
#      B1: "Unemployment %(y)"   -> "Unemployment %"
#  - The "Data" sheet becomes the active sheet/tab (it was "Description"),
#    with the on-sheet selection moved to D5.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")

$ws1.Range("A1").Value = "Interest Rates"
$ws1.Range("B1").Value = "Unemployment %"

# Make "Data" the active/selected sheet and move the selection on it.
$ws1.Activate() | Out-Null
$ws1.Range("D5").Select() | Out-Null
